# Update cohort and re-run modules: refresh p-values produced by the
# demo-clin-stats pipeline for each tumor-histology worksheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    "Low-grade glioma" = @{ "C3" = 0.600739926007399; "C6" = 0.0693930606939306; "C7" = 0.226177382261774; "C8" = 0.772022797720228 }
    "Non-neoplastic tumor" = @{ "C3" = 0.912108789121088; "C5" = 0.0001999800019998; "C6" = 0.253774622537746; "C7" = 0.851214878512149 }
    "Mixed neuronal-glial tumor" = @{ "C3" = 0.188581141885811; "C6" = 0.515048495150485; "C7" = 0.787321267873213; "C8" = 0.693730626937306 }
    "Medulloblastoma" = @{ "C3" = 0.776522347765223; "C6" = 0.482451754824518; "C7" = 0.296670332966703; "C8" = 0.837416258374163 }
    "Schwannoma" = @{ "C2" = 0.367226758676688; "C3" = 0.802919708029197; "C5" = 0.499150084991501; "C6" = 0.928107189281072; "C7" = 0.84021597840216 }
    "Mesenchymal tumor" = @{ "C3" = 0.127987201279872; "C5" = 0.0004999500049995; "C6" = 0.0233976602339766; "C7" = 0.842915708429157; "C8" = 0.491050894910509 }
    "Germ cell tumor" = @{ "C3" = 0.614438556144386; "C4" = 0.0001999800019998; "C5" = 0.0094990500949905; "C6" = 0.856814318568143; "C7" = 0.368763123687631; "C8" = 0.286371362863714 }
    "Craniopharyngioma" = @{ "C3" = 0.742225777422258; "C5" = 0.0053994600539946; "C6" = 0.167883211678832; "C7" = 0.582841715828417 }
    "Other tumor" = @{ "C2" = 0.660987736725295; "C3" = 0.0426957304269573; "C5" = 0.0011998800119988; "C6" = 0.877512248775122; "C7" = 0.993600639936006; "C8" = 0.555138861058704; "C9" = 0.655966505077126 }
    "Ependymoma" = @{ "C3" = 0.398160183981602; "C6" = 0.261273872612739; "C7" = 0.278572142785721; "C8" = 0.965103489651035 }
    "DIPG or DMG" = @{ "C3" = 0.501649835016498; "C6" = 0.014998500149985; "C7" = 0.617738226177382 }
    "ATRT" = @{ "C3" = 0.367463253674633; "C6" = 0.494850514948505; "C7" = 0.276172382761724; "C8" = 0.173582641735826 }
    "Other high-grade glioma" = @{ "C3" = 0.800919908009199; "C6" = 0.84011598840116; "C7" = 0.0145985401459854; "C8" = 0.667333266673333 }
    "Meningioma" = @{ "C3" = 0.944305569443056; "C5" = 0.0003999600039996; "C6" = 0.677232276772323; "C7" = 0.30976902309769 }
    "Neurofibroma plexiform" = @{ "C4" = 0.0002999700029997; "C5" = 0.149385061493851; "C6" = 0.0606939306069393; "C7" = 0.301569843015698 }
    "Oligodendroglioma" = @{ "C3" = 0.163383661633837; "C5" = 0.209379062093791; "C7" = 0.211578842115788 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
